# Auto-generated: apply cryptos list price/volume refresh (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.458.01"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.683.01"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'685.97"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'159.86"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'7.04"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").Value = "'0.434"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "4.304.47"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'32.36"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "69.424.84"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.665.40"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("D20").Value = "'471.55"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "'9.95"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "'79.62"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "3.831.51"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "'0.0000125"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").Value = "'9.21"
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -5.63%  "
$ws.Range("D32").Value = "'6.57"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "'26.88"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "3.657.73"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").Value = "'8.20"
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("D38").Value = "'6.14"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'2.23"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "'0.0898"
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'165.74"
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'47.57"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "'1.12"
$ws.Range("E48").Value = "  +5.59%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'27.68"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "'7.77"
$ws.Range("E51").Value = "  -3.28%  "
